# Profile and ModelCode update
#
# 1) Sheet1 ("abstract-concrete"): add a new column A listing the CIM
#    super-class for each concrete class in column B, plus a new row 8
#    for "Wires::Switch"; also fix a typo in an existing shared string
#    used by B7 (Core::RegularIntervalTimePoint -> Core::RegularTimePoint).
# 2) Sheet2: renamed to "attributes" and populated with a 7-column
#    attribute table; it becomes the active/selected sheet.
#
# NOTE on write order: the shared-string table is appended to in the
# exact order new distinct strings are first written, and a string that
# becomes orphaned (its only referencing cell is overwritten) is dropped
# from the table immediately, before later appends -- so the order of
# the statements below is chosen deliberately to reproduce the exact
# shared-string layout of the target workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: abstract-concrete
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Fix the existing shared string referenced by B7 first, while it is
# still the only/last entry in the table, so the corrected text reuses
# slot 7 instead of being appended after the new strings below.
$ws1.Range("B7").Value = "Core::RegularTimePoint"

$ws1.Range("A2").Value = "Core::IdentifiedObject"
$ws1.Range("A3").Value = "Core::PowerSystemResource"
$ws1.Range("A4").Value = "Core::Equipment"
$ws1.Range("A5").Value = "Core::ConductingEquipment"

# Row 8 (Wires::Switch) is written before rows 6/7 so the new shared
# string for "Wires::Switch" lands ahead of BasicIntervalSchedule /
# IrregularIntervalSchedule, matching the source ordering.
$ws1.Range("A8").Value = "Wires::Switch"

$ws1.Range("A6").Value = "Core::BasicIntervalSchedule"
$ws1.Range("A7").Value = "Core::IrregularIntervalSchedule"

# ---------------------------------------------------------------------
# Sheet2: renamed to "attributes"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "attributes"

# Column widths (best effort -- the COM layer quantizes to 1/6 of a
# character, so these are the closest achievable values).
$ws2.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws2.Columns.Item(2).ColumnWidth = 28.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 26.333333333333332
$ws2.Columns.Item(4).ColumnWidth = 25.666666666666668
$ws2.Columns.Item(5).ColumnWidth = 21.666666666666668
$ws2.Columns.Item(6).ColumnWidth = 25.666666666666668
$ws2.Columns.Item(7).ColumnWidth = 29.166666666666668

# Populate column-by-column (not row-by-row) so new shared strings are
# created in the same order as the target file.

# Column A
$ws2.Range("A1").Value = "IdentifiedObject"
$ws2.Range("A2").Value = "aliasName: String"
$ws2.Range("A3").Value = "mRID: String"
$ws2.Range("A4").Value = "name: String"

# Column B
$ws2.Range("B1").Value = "BasicIntervalSchedule"
$ws2.Range("B2").Value = "startTime: DateTime"
$ws2.Range("B3").Value = "value1Multiplier: UnitMultiplier"
$ws2.Range("B4").Value = "value1Unit: UnitSymbol"
$ws2.Range("B5").Value = "value2Multiplier: UnitMultiplier"
$ws2.Range("B6").Value = "value2Unit: UnitSymbol"

# Column C
$ws2.Range("C1").Value = "SwitchingOperation"
$ws2.Range("C2").Value = "newState: SwitchState"
$ws2.Range("C3").Value = "operationTime: DateTime"
$ws2.Range("C4").Value = "OutageSchedule: REFERENCE"

# Column D
$ws2.Range("D1").Value = "IrregularTimePoint"
$ws2.Range("D2").Value = "time: Secounds"
$ws2.Range("D3").Value = "value1: Float"
$ws2.Range("D4").Value = "value2: Float"
$ws2.Range("D5").Value = "IntervalSchedule: REFERENCE"

# Column E
$ws2.Range("E1").Value = "RegularIntervalSchedule"
$ws2.Range("E2").Value = "endTime: DateTime"
$ws2.Range("E3").Value = "timeStep: Secounds"

# Column F (F3:F5 reuse the D3:D5 strings created above)
$ws2.Range("F1").Value = "RegularTimePoint"
$ws2.Range("F2").Value = "sequenceNumber: Integer"
$ws2.Range("F3").Value = "value1: Float"
$ws2.Range("F4").Value = "value2: Float"
$ws2.Range("F5").Value = "IntervalSchedule: REFERENCE"

# Column G
$ws2.Range("G1").Value = "Switch"
$ws2.Range("G2").Value = "SwitchingOperations: REFERENCE"

# ---------------------------------------------------------------------
# Selection / active-sheet state
# ---------------------------------------------------------------------
$ws1.Range("A8").Select()
$ws2.Activate()
$ws2.Range("G8").Select()
